# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.287.13"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "1.869.71"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'235.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").Value = "'0.4700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.72%  "

# Row 8
$ws.Range("D8").Value = "'0.2875"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

# Row 9
$ws.Range("D9").Value = "'0.06582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.06%  "

# Row 10
$ws.Range("D10").Value = "'21.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "

# Row 11
$ws.Range("D11").Value = "'0.08020"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "

# Row 12
$ws.Range("D12").Value = "'97.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").Value = "1.874.56"
$ws.Range("E13").Value = "  +0.80%  "

# Row 14
$ws.Range("D14").Value = "'5.129"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "

# Row 15
$ws.Range("D15").Value = "'0.6856"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "

# Row 16
$ws.Range("D16").Value = "'269.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.32%  "

# Row 17
$ws.Range("D17").Value = "30.273.39"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18
$ws.Range("D18").Value = "'14.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19
$ws.Range("D19").Value = "'0.000007670"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.06%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").Value = "2.117.64"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").Value = "'5.276"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.63%  "

# Row 24
$ws.Range("D24").Value = "'6.219"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'9.431"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "

# Row 26
$ws.Range("E26").Value = "  +1.10%  "

# Row 27
$ws.Range("D27").Value = "'18.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "

# Row 28
$ws.Range("E28").Value = "  +1.52%  "

# Row 29
$ws.Range("D29").Value = "'1.369"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

# Row 30
$ws.Range("D30").Value = "'0.09884"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.01%  "

# Row 31
$ws.Range("D31").Value = "'4.380"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("D32").Value = "'1.463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "

# Row 33
$ws.Range("D33").Value = "'4.077"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("D34").Value = "'0.04706"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "

# Row 35
$ws.Range("D35").Value = "'1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "

# Row 36
$ws.Range("D36").Value = "'0.7007"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "

# Row 37
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("D38").Value = "'0.01876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "

# Row 39
$ws.Range("E39").Value = "  +0.26%  "

# Row 40
$ws.Range("D40").Value = "'6.297"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("D41").Value = "'72.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.66%  "

# Row 42
$ws.Range("D42").Value = "'1.953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.49%  "

# Row 43
$ws.Range("D43").Value = "'0.8431"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4166"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").Value = "'103.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("E47").Value = "  -0.19%  "

# Row 48
$ws.Range("D48").Value = "'7.061"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
$ws.Range("D49").Value = "'924.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.87%  "

# Row 50
$ws.Range("D50").Value = "'34.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.70%  "

# Row 51
$ws.Range("D51").Value = "'0.05675"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
